# Fruta / hortaliza, semanal
# Insert 3 new daily price rows for "Uva" (grape) right before the existing
# row 761 block, which pushes the old rows 761-816 down to 764-819 and grows
# the used range from A1:T816 to A1:T819.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Push rows 761-816 down by 3 rows (Excel copies formatting, incl. the date
# number format on column D, from the row immediately above the insertion).
$ws.Rows("761:763").Insert()

# Row 761 - new "Autumn Royal" entry
$ws.Range("A761").Value = 9
$ws.Range("B761").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C761").Value = "Metropolitana"
$ws.Range("D761").Value = 45013
$ws.Range("E761").Value = 13
$ws.Range("F761").Value = "Fruta"
$ws.Range("G761").Value = 100109
$ws.Range("H761").Value = "Uva"
$ws.Range("I761").Value = 100109001
$ws.Range("J761").Value = "Uva"
$ws.Range("K761").Value = "Autumn Royal"
$ws.Range("L761").Value = "Primera"
$ws.Range("M761").Value = 200
$ws.Range("N761").Value = 10000
$ws.Range("O761").Value = 10000
$ws.Range("P761").Value = 10000
$ws.Range("Q761").Value = "$/bandeja 18 kilos"
$ws.Range("R761").Value = "Región de O'Higgins"
$ws.Range("S761").Value = 556
$ws.Range("T761").Value = 18

# Row 762 - new "Red Globe" entry
$ws.Range("A762").Value = 9
$ws.Range("B762").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C762").Value = "Metropolitana"
$ws.Range("D762").Value = 45013
$ws.Range("E762").Value = 13
$ws.Range("F762").Value = "Fruta"
$ws.Range("G762").Value = 100109
$ws.Range("H762").Value = "Uva"
$ws.Range("I762").Value = 100109001
$ws.Range("J762").Value = "Uva"
$ws.Range("K762").Value = "Red Globe"
$ws.Range("L762").Value = "Primera"
$ws.Range("M762").Value = 300
$ws.Range("N762").Value = 10000
$ws.Range("O762").Value = 10000
$ws.Range("P762").Value = 10000
$ws.Range("Q762").Value = "$/bandeja 18 kilos"
$ws.Range("R762").Value = "Provincia de Curicó"
$ws.Range("S762").Value = 556
$ws.Range("T762").Value = 18

# Row 763 - new "Thompson seedless" entry
$ws.Range("A763").Value = 9
$ws.Range("B763").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C763").Value = "Metropolitana"
$ws.Range("D763").Value = 45013
$ws.Range("E763").Value = 13
$ws.Range("F763").Value = "Fruta"
$ws.Range("G763").Value = 100109
$ws.Range("H763").Value = "Uva"
$ws.Range("I763").Value = 100109001
$ws.Range("J763").Value = "Uva"
$ws.Range("K763").Value = "Thompson seedless"
$ws.Range("L763").Value = "Primera"
$ws.Range("M763").Value = 300
$ws.Range("N763").Value = 11000
$ws.Range("O763").Value = 11000
$ws.Range("P763").Value = 11000
$ws.Range("Q763").Value = "$/bandeja 18 kilos"
$ws.Range("R763").Value = "Región de O'Higgins"
$ws.Range("S763").Value = 611
$ws.Range("T763").Value = 18
